$wb = $excel.ActiveWorkbook

# --- Sheet "Averages": remove "PyOMP with JIT" column, add new Pi Computation data ---
$wsAvg = $wb.Worksheets.Item("Averages")
$wsAvg.Columns.Item(5).Delete()  # delete column E ("PyOMP with JIT")

# Pi Computation row is row 3; PyOMP column is D; new Serial (Python) col is F
$wsAvg.Range("D3").Value = 0.016643667
$wsAvg.Range("F3").Value = 5.869972436

# --- Sheet "Pi Computation C++": add PyOMP (16 threads) header first ---
$wsPi = $wb.Worksheets.Item("Pi Computation C++")
$wsPi.Range("C1").Value = "PyOMP (16 threads)"

# --- Sheet "Information": add note cell C2 ---
$wsInfo = $wb.Worksheets.Item("Information")
$wsInfo.Range("C2").Value = "The PyOMP time values are without the JIT compilation time"

# --- Sheet "Pi Computation C++": add Python (Serial) column header + data ---
$wsPi.Range("D1").Value = "Python (Serial)"

$cVals = @(0.0195701122283935, 0.014625072479248, 0.0244920253753662, 0.0146269798278808, 0.0165870189666748, 0.0144770145416259, 0.0153179168701171, 0.0169758796691894, 0.0176389217376709, 0.0169880390167236, 0.0157358646392822, 0.0164859294891357, 0.0154941082000732, 0.0154900550842285, 0.0151500701904296)
$dVals = @(5.84530711174011, 5.84022879600524, 5.89002227783203, 5.85703229904174, 5.78420734405517, 5.83787512779235, 5.83210206031799, 6.06160354614257, 5.79919981956481, 5.83363389968872, 5.88991665840148, 5.80106353759765, 5.82321572303772, 5.95681309700012, 5.99736523628234)

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 2
    $wsPi.Range("C$row").Value = $cVals[$i]
    $wsPi.Range("D$row").Value = $dVals[$i]
}

# --- Selections / active sheet ---
$wsInfo.Range("C17").Select()
$wsAvg.Range("F5").Select()
$wsPi.Activate()
$wsPi.Range("D16").Select()
